# Apply "Hjemme passive tweaks lichtwark deleted values" edit:
# updates the first four data columns (B:E) of rows 1-3 on the only
# worksheet, and narrows the active selection from B1:AY3 to B1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header / count values)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON)
$ws.Range("B2").Value = 24.832978138124986
$ws.Range("C2").Value = 23.28494295750005
$ws.Range("D2").Value = 31.556935828125006
$ws.Range("E2").Value = 34.665075633125014

# Row 3 (STR)
$ws.Range("B3").Value = 20.167593688124953
$ws.Range("C3").Value = 38.205168322500015
$ws.Range("D3").Value = 35.443980539999927
$ws.Range("E3").Value = 31.535374552500002

# Shrink the saved selection/active range to match the new data extent
$ws.Range("B1:E3").Select() | Out-Null
